$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: split the "Pupil area = ..." run into two runs:
#   "Pupil area = "  +  "Complexity Index (derived from refined multi-scale entropy)"
# ------------------------------------------------------------------

# Shield the future run boundary ("Pupil area = ") with a transient
# format toggle so the replace below won't let it re-merge with the
# text that follows it.
$r1 = $d.Content
$r1.Find.Execute("Pupil area = ")
$boundary1 = $d.Range($r1.Start, $r1.End)
$boundary1.Bold = 1

# Replace the remainder of the sentence with the new wording.
$r2 = $d.Content
$r2.Find.Execute("inverse memory level index (derived from Hurst exponents)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Complexity Index (derived from refined multi-scale entropy)", 2)

# Remove the transient bold from the "Pupil area = " run again.
$r3 = $d.Content
$r3.Find.Execute("Pupil area = ")
$boundary1b = $d.Range($r3.Start, $r3.End)
$boundary1b.Bold = 0

# Materialize an explicit (empty) run-properties element on the new
# "Complexity Index ..." run to match the document's normal style.
$r4 = $d.Content
$r4.Find.Execute("Complexity Index (derived from refined multi-scale entropy)")
$newRun = $d.Range($r4.Start, $r4.End)
$newRun.Bold = 1
$newRun.Bold = 0

# ------------------------------------------------------------------
# Change 2: merge the two "Ear width ..." runs into a single run:
#   "Ear width = evenness or control of beat intervals " + "(derived from sums of square differences)"
#   -> "Ear width = evenness or control of beat intervals (derived from sums of square differences)"
# ------------------------------------------------------------------

# Shield the preceding leading-whitespace run so it will not be pulled
# into the merge triggered by the replace below.
$e1 = $d.Content
$e1.Find.Execute("Ear width")
$leading = $d.Range($e1.Start - 4, $e1.Start)
$leading.Bold = 1

# Perform a no-text-change "replace" across both runs: this causes the
# COM layer to normalize/merge adjacent same-formatted runs.
$e2 = $d.Content
$e2.Find.Execute("Ear width = evenness or control of beat intervals (derived from sums of square differences)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Ear width = evenness or control of beat intervals (derived from sums of square differences)", 2)

# Remove the transient bold from the leading-whitespace run again.
$e3 = $d.Content
$e3.Find.Execute("Ear width")
$leadingB = $d.Range($e3.Start - 4, $e3.Start)
$leadingB.Bold = 0

# Materialize an explicit (empty) run-properties element on the merged
# run to match the document's normal style.
$e4 = $d.Content
$e4.Find.Execute("Ear width = evenness or control of beat intervals (derived from sums of square differences)")
$mergedRun = $d.Range($e4.Start, $e4.End)
$mergedRun.Bold = 1
$mergedRun.Bold = 0
